$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 131 (shifts existing rows 131-150 down to 132-151)
$ws.Rows(131).Insert()

# Populate the newly inserted row 131 with the new weekly data point
$ws.Range("A131").Value = 3
$ws.Range("B131").Value = "Femacal de La Calera"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 44522
$ws.Range("E131").Value = 5
$ws.Range("F131").Value = 100112010
$ws.Range("G131").Value = "Achicoria"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 60
$ws.Range("K131").Value = 6000
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = 6000
$ws.Range("N131").Value = "$/caja 16 unidades"
$ws.Range("O131").Value = "Provincia de Quillota"
$ws.Range("P131").Value = 375
$ws.Range("Q131").Value = 16
$ws.Range("R131").Value = "Hortaliza"
